$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (scheduled GitHub Actions update).
# Column D ("Price") values such as "603.90", "31.60" or "1.00" look like
# plain numbers to Excel and would otherwise get silently reinterpreted as
# numeric values, losing the exact text formatting used by the source data
# (e.g. trailing zeros, or thousand-dot grouped strings like "66.701.32").
# Force those cells to Text format before writing the new value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '66.701.32'
$ws.Range("E2").Value = '  -2.54%  '

# Row 3
$ws.Range("D3").Value = '3.479.50'
$ws.Range("E3").Value = '  -2.19%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '603.90'
$ws.Range("E5").Value = '  -2.58%  '

# Row 6
$ws.Range("D6").Value = '148.25'
$ws.Range("E6").Value = '  -4.79%  '

# Row 7
$ws.Range("D7").Value = '3.477.05'
$ws.Range("E7").Value = '  -2.27%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").Value = '0.482'
$ws.Range("E9").Value = '  -1.54%  '

# Row 10
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  -3.22%  '

# Row 11
$ws.Range("D11").Value = '7.57'
$ws.Range("E11").Value = '  +3.03%  '

# Row 12
$ws.Range("D12").Value = '0.425'
$ws.Range("E12").Value = '  -3.48%  '

# Row 13
$ws.Range("D13").Value = '0.0000214'
$ws.Range("E13").Value = '  -4.74%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '4.068.78'
$ws.Range("E14").Value = '  -2.25%  '

# Row 15
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '31.60'
$ws.Range("E15").Value = '  -5.08%  '

# Row 16
$ws.Range("D16").Value = '3.494.77'
$ws.Range("E16").Value = '  -1.86%  '

# Row 17
$ws.Range("D17").Value = '66.831.74'
$ws.Range("E17").Value = '  -2.76%  '

# Row 18
$ws.Range("E18").Value = '  -0.45%  '

# Row 19
$ws.Range("D19").Value = '6.47'
$ws.Range("E19").Value = '  -4.68%  '

# Row 20
$ws.Range("D20").Value = '15.40'
$ws.Range("E20").Value = '  -3.73%  '

# Row 21
$ws.Range("D21").Value = '10.10'
$ws.Range("E21").Value = '  +0.34%  '

# Row 22
$ws.Range("D22").Value = '439.51'
$ws.Range("E22").Value = '  -4.25%  '

# Row 23
$ws.Range("D23").Value = '0.611'
$ws.Range("E23").Value = '  -5.02%  '

# Row 24
$ws.Range("D24").Value = '79.62'
$ws.Range("E24").Value = '  +1.28%  '

# Row 25
$ws.Range("E25").Value = '  +0.05%  '

# Row 26
$ws.Range("D26").Value = '3.616.47'
$ws.Range("E26").Value = '  -2.39%  '

# Row 27
$ws.Range("D27").Value = '0.0000120'
$ws.Range("E27").Value = '  -8.51%  '

# Row 28
$ws.Range("D28").Value = '9.76'
$ws.Range("E28").Value = '  -7.52%  '

# Row 29
$ws.Range("D29").Value = '8.39'
$ws.Range("E29").Value = '  -7.44%  '

# Row 30
$ws.Range("D30").Value = '2.49'
$ws.Range("E30").Value = '  -2.83%  '

# Row 31
$ws.Range("D31").Value = '1.59'
$ws.Range("E31").Value = '  -6.11%  '

# Row 32
$ws.Range("E32").Value = '  -0.97%  '

# Row 33
$ws.Range("E33").Value = '  +0.04%  '

# Row 34
$ws.Range("D34").Value = '25.42'
$ws.Range("E34").Value = '  -3.24%  '

# Row 35
$ws.Range("D35").Value = '6.05'
$ws.Range("E35").Value = '  -6.47%  '

# Row 36
$ws.Range("D36").Value = '3.469.58'
$ws.Range("E36").Value = '  -2.36%  '

# Row 37
$ws.Range("D37").Value = '1.80'
$ws.Range("E37").Value = '  -6.76%  '

# Row 38
$ws.Range("D38").Value = '7.93'
$ws.Range("E38").Value = '  -4.51%  '

# Row 39
$ws.Range("E39").Value = '  +0.00%  '

# Row 40
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.00%  '

# Row 41
$ws.Range("D41").Value = '175.95'
$ws.Range("E41").Value = '  -1.50%  '

# Row 42
$ws.Range("D42").Value = '0.0889'
$ws.Range("E42").Value = '  -3.23%  '

# Row 43
$ws.Range("E43").Value = '  -10.43%  '

# Row 44
$ws.Range("D44").Value = '5.42'
$ws.Range("E44").Value = '  -3.31%  '

# Row 45
$ws.Range("E45").Value = '  -1.17%  '

# Row 46
$ws.Range("D46").Value = '28.95'
$ws.Range("E46").Value = '  -5.26%  '

# Row 47
$ws.Range("D47").Value = '46.27'
$ws.Range("E47").Value = '  +1.05%  '

# Row 48
$ws.Range("E48").Value = '  -8.00%  '

# Row 49
$ws.Range("D49").Value = '7.47'
$ws.Range("E49").Value = '  -4.29%  '

# Row 50
$ws.Range("D50").Value = '2.43'
$ws.Range("E50").Value = '  -9.26%  '

# Row 51
$ws.Range("D51").Value = '0.983'
$ws.Range("E51").Value = '  -4.29%  '

# Restore the default ("Normal") style on the Price cells we just force-
# formatted as Text, so their appearance matches the rest of the sheet -
# only the underlying text type is preserved, not the Text number format.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
